$d = $word.ActiveDocument

# The table's 2nd and 3rd columns ("Unique (#)" / "Missing (%)" header,
# with their corresponding data cells in every row) are being dropped
# entirely, shrinking the table from 8 to 6 columns.
$t = $d.Tables.Item(1)
$t.Columns.Item(2).Delete()
$t.Columns.Item(2).Delete()

# Rename the row-label cells to the more descriptive labels used by the
# updated modelsummary output.
$d.Content.Find.Execute("schooling", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Schooling (Yrs)", 2)
$d.Content.Find.Execute("life_expectancy", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Life Expectancy (Yrs)", 2)
